$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "355.45", "0.570", "52.031.00") are preserved verbatim as text,
# matching the source inlineStr cells instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '52.031.00'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '2.819.77'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '355.45'
$ws.Range("E5").Value = '  +2.74%  '
$ws.Range("D6").Value = '111.94'
$ws.Range("E6").Value = '  -4.17%  '
$ws.Range("D7").Value = '0.570'
$ws.Range("E7").Value = '  +3.04%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.598'
$ws.Range("E9").Value = '  +2.23%  '
$ws.Range("D10").Value = '40.72'
$ws.Range("E10").Value = '  -6.27%  '
$ws.Range("D11").Value = '0.0860'
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").Value = '19.92'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("D14").Value = '7.75'
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").Value = '3.259.58'
$ws.Range("E15").Value = '  +0.75%  '
$ws.Range("D16").Value = '2.819.54'
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("D17").Value = '0.931'
$ws.Range("E17").Value = '  +4.27%  '
$ws.Range("D18").Value = '51.802.10'
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("D19").Value = '7.53'
$ws.Range("E19").Value = '  +5.59%  '
$ws.Range("D20").Value = '3.19'
$ws.Range("E20").Value = '  -0.75%  '
$ws.Range("D21").Value = '13.43'
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = '0.0₃0995'
$ws.Range("E22").Value = '  +1.26%  '
$ws.Range("D23").Value = '70.80'
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("D24").Value = '269.22'
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("D25").Value = '2.82'
$ws.Range("E25").Value = '  +2.16%  '
$ws.Range("D26").Value = '26.98'
$ws.Range("E26").Value = '  +1.13%  '
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").Value = '10.32'
$ws.Range("E28").Value = '  +0.53%  '
$ws.Range("D29").Value = '2.26'
$ws.Range("E29").Value = '  +1.01%  '
$ws.Range("D30").Value = '0.0489'
$ws.Range("E30").Value = '  +17.98%  '
$ws.Range("E31").Value = '  +2.38%  '
$ws.Range("D32").Value = '52.50'
$ws.Range("E32").Value = '  +4.42%  '
$ws.Range("D33").Value = '34.79'
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("D34").Value = '5.95'
$ws.Range("E34").Value = '  +3.76%  '
$ws.Range("D35").Value = '5.60'
$ws.Range("E35").Value = '  +12.35%  '
$ws.Range("D36").Value = '0.0852'
$ws.Range("E36").Value = '  +3.51%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  +1.33%  '
$ws.Range("E39").Value = '  -4.12%  '
$ws.Range("D40").Value = '18.37'
$ws.Range("E40").Value = '  -3.11%  '
$ws.Range("E41").Value = '  +1.09%  '
$ws.Range("D42").Value = '126.85'
$ws.Range("E42").Value = '  -0.96%  '
$ws.Range("D43").Value = '23.13'
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("E44").Value = '  -8.02%  '
$ws.Range("E45").Value = '  -1.63%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '3.35'
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.078.41'
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("E48").Value = '  -4.84%  '
$ws.Range("D49").Value = '5.87'
$ws.Range("E49").Value = '  +6.05%  '
$ws.Range("D50").Value = '0.975'
$ws.Range("E50").Value = '  -0.81%  '
$ws.Range("D51").Value = '9.14'
$ws.Range("E51").Value = '  +1.97%  '
